$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.304.25'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.747.36'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.24'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.57'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.746.86'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.541'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.172'
$ws.Range("E10").Value = '  +4.89%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.27'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.376.00'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.744.35'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.262.27'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.42'
$ws.Range("E18").Value = '  +2.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.42'
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.31'
$ws.Range("E21").Value = '  +12.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.29'
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  +5.22%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.32'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.98'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.64'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.895.33'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.682.20'
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  +6.09%  '
$ws.Range("E39").Value = '  +3.20%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.327'
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("E42").Value = '  +6.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.89'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '423.81'
$ws.Range("E45").Value = '  -2.30%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.17'
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.19'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.794.85'
$ws.Range("E50").Value = '  +1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0355'
$ws.Range("E51").Value = '  +0.76%  '
